$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update events-toci (F) and events-control (H) values for the RECOVERY trial rows (17-26)
# after RECOVERY was published on Lancet.

$ws.Range("F17").Value = 180
$ws.Range("H17").Value = 214

$ws.Range("F18").Value = 310
$ws.Range("H18").Value = 366

$ws.Range("F19").Value = 131
$ws.Range("H19").Value = 149

$ws.Range("F20").Value = 697
$ws.Range("H20").Value = 635

$ws.Range("F21").Value = 401
$ws.Range("H21").Value = 362

$ws.Range("F22").Value = 52

$ws.Range("F23").Value = 214
$ws.Range("H23").Value = 256

$ws.Range("F24").Value = 407
$ws.Range("H24").Value = 473

$ws.Range("F25").Value = 380
$ws.Range("H25").Value = 331

$ws.Range("F26").Value = 770
$ws.Range("H26").Value = 713

# Update the view state: select E20 (matches the saved selection in the workbook)
$ws.Activate()
$ws.Range("E20").Select()
